$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The table currently occupies rows 2-18 (header row 1, data rows 2-7, then
# empty formatted rows 8-18). We need to:
#   1) Insert a new data row before the current row 7 (shifting the existing
#      row 7 "Portfolio / Self" entry down to row 8, and all the blank rows
#      down by one as well).
#   2) Add one extra blank row at the bottom (row 20) so the blank-row pool
#      stays the same size as before.
#   3) Fill the new row 7 with the "Stand Up recording" feedback entry.
#   4) Update the sheet view (top-left cell / selection).
# We shift content with Range.Copy so the existing cell styles (border,
# number format, wrap, vertical alignment) are reused instead of new style
# entries being synthesized.
# ---------------------------------------------------------------------------

# Shift rows 18 down through 7 down by one row (process bottom-up so we never
# overwrite a row before it has been copied away).
for ($r = 18; $r -ge 7; $r--) {
    $src = $ws.Range("A" + $r + ":F" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":F" + ($r + 1))
    $src.Copy($dst)
}

# Row 8 is now a copy of the old row 7 ("Portfolio" / "Self" / developing
# portfolio). Restore its custom row height (lost because only the A:F cell
# range was copied, not the row height).
$ws.Rows.Item(8).RowHeight = 43.5

# Row 7 still holds a duplicate of the old row 7 content/style at this point
# (it was only copied away, not cleared). Reuse that formatting, just bump
# its height to match the new feedback's wrapped text (58, matching rows
# 3/4 which hold similarly long text) and overwrite its values below.
$ws.Rows.Item(7).RowHeight = 58

# New row 20: copy formatting from the (still blank) row 19 so it matches
# the other blank rows.
$ws.Range("A19:F19").Copy($ws.Range("A20:F20"))

# ---------------------------------------------------------------------------
# Populate the newly inserted row 7 with the new feedback entry.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 44977
$ws.Range("B7").Value = "Portfolio"
$ws.Range("C7").Value = "Facilitator"
$ws.Range("D7").Value = "Stand Up recording, as advised by Module 3.4 is not required for PortfolioRecord this stand-up and add it to your portfolio"
$ws.Range("E7").Value = "Removed the stand Up Recording from Portfolio"
$ws.Range("F7").Value = "Completed"

# ---------------------------------------------------------------------------
# Update the sheet view to match the saved selection/scroll position.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("A8").Select()
